$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely; row 3 shifts up to become the new row 2.
$ws.Rows.Item(2).Delete()
